# Rename the "hiddenfigures" defined names (one per sheet scope) to "sampha1".
# NOTE: the Names collection is "live" by index, so renaming front-to-back
# causes Item() to skip entries after a rename shifts things around;
# iterating back-to-front keeps each Item($i) pointing at the name we
# actually intend to touch.
$wb = $excel.ActiveWorkbook
$namesCount = $wb.Names.Count
for ($i = $namesCount; $i -ge 1; $i--) {
    $wb.Names.Item($i).Name = "sampha1"
}

# New tracklist data (Sampha's "Process" bootlegs/old reviews) replacing the
# previous Pharrell-related tracklist, for both Sheet1 and Sheet3 (which
# mirror each other - both are driven by the same web query).
$titles = @(
    "Plastic 100°C",
    "Blood on Me",
    "Kora Sings",
    "(No One Knows Me) Like the Piano",
    "Take Me Inside",
    "Reverse Faults",
    "Under",
    "Timmy's Prayer",
    "Incomplete Kisses",
    "What Shouldn't I Be?"
)
$composers = @(
    "Sampha Sisay",
    "Sampha Sisay",
    "Sampha Sisay",
    "Sampha Sisay",
    "Sampha Sisay",
    "Sampha Sisay",
    "Sampha Sisay",
    "Sampha Sisay / Kanye West",
    "Sampha Sisay",
    "Sampha Sisay"
)
$performers = @(
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha",
    "Sampha"
)
$times = @(
    0.21944444444444444,
    0.17083333333333331,
    0.17847222222222223,
    0.15138888888888888,
    0.096527777777777768,
    0.17569444444444446,
    0.19513888888888889,
    0.18263888888888891,
    0.16180555555555556,
    0.14722222222222223
)

foreach ($sheetName in @("Sheet1", "Sheet3")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($r = 0; $r -lt 10; $r++) {
        $row = $r + 2
        $ws.Range("B$row").Value = $titles[$r]
        $ws.Range("C$row").Value = $composers[$r]
        $ws.Range("D$row").Value = $performers[$r]
        $ws.Range("E$row").Value = $times[$r]
    }

    # Column widths: B grows to fit the longer titles, C shrinks slightly,
    # D shrinks a lot now that "Sampha" is the only performer.
    # (ColumnWidth is expressed in characters; the host rounds to whole
    # pixels at save time, so we pick the character width whose rounded
    # pixel-width lands on the desired OOXML column width.)
    $ws.Columns.Item(2).ColumnWidth = 31.93
    $ws.Columns.Item(3).ColumnWidth = 25.65
    $ws.Columns.Item(4).ColumnWidth = 9.2857142857142857
}

# Sheet1's selection moved from A1:E11 (active cell E11) to just B27.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("B27").Select()
